$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1431
$ws.Range("B3").Value = 2954
$ws.Range("B4").Value = 514
$ws.Range("B5").Value = 1167
$ws.Range("B6").Value = 1273
$ws.Range("B7").Value = 3052
$ws.Range("B8").Value = 257
$ws.Range("B9").Value = 332
$ws.Range("B10").Value = 2463
